$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "57.946.18"
$ws.Cells.Item(2, 5).Value = "  +0.35%  "
$ws.Cells.Item(3, 4).Value = "2.344.87"
$ws.Cells.Item(3, 5).Value = "  +0.12%  "
$ws.Cells.Item(4, 5).Value = "  +0.14%  "
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = "540.72"
$ws.Cells.Item(5, 4).Style = "Normal"
$ws.Cells.Item(5, 5).Value = "  -0.27%  "
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = "134.09"
$ws.Cells.Item(6, 4).Style = "Normal"
$ws.Cells.Item(6, 5).Value = "  -0.64%  "
$ws.Cells.Item(7, 5).Value = "  +0.08%  "
$ws.Cells.Item(8, 4).NumberFormat = "@"
$ws.Cells.Item(8, 4).Value = "0.562"
$ws.Cells.Item(8, 4).Style = "Normal"
$ws.Cells.Item(8, 5).Value = "  +4.37%  "
$ws.Cells.Item(9, 5).Value = "  +0.54%  "
$ws.Cells.Item(10, 5).Value = "  +1.82%  "
$ws.Cells.Item(11, 5).Value = "  -1.83%  "
$ws.Cells.Item(12, 5).Value = "  +0.21%  "
$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = "23.81"
$ws.Cells.Item(13, 4).Style = "Normal"
$ws.Cells.Item(13, 5).Value = "  +1.07%  "
$ws.Cells.Item(14, 4).Value = "2.762.16"
$ws.Cells.Item(14, 5).Value = "  -0.72%  "
$ws.Cells.Item(15, 4).Value = "57.896.12"
$ws.Cells.Item(15, 5).Value = "  +0.15%  "
$ws.Cells.Item(16, 5).Value = "  +0.36%  "
$ws.Cells.Item(17, 4).Value = "2.344.75"
$ws.Cells.Item(17, 5).Value = "  -1.52%  "
$ws.Cells.Item(18, 5).Value = "  +0.76%  "
$ws.Cells.Item(19, 5).Value = "  +2.04%  "
$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = "328.69"
$ws.Cells.Item(20, 4).Style = "Normal"
$ws.Cells.Item(20, 5).Value = "  -1.85%  "
$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = "6.73"
$ws.Cells.Item(21, 4).Style = "Normal"
$ws.Cells.Item(21, 5).Value = "  -0.66%  "
$ws.Cells.Item(22, 5).Value = "  -0.03%  "
$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = "62.94"
$ws.Cells.Item(23, 4).Style = "Normal"
$ws.Cells.Item(23, 5).Value = "  +1.31%  "
$ws.Cells.Item(24, 5).Value = "  -3.38%  "
$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = "0.993"
$ws.Cells.Item(25, 4).Style = "Normal"
$ws.Cells.Item(25, 5).Value = "  -0.38%  "
$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = "8.30"
$ws.Cells.Item(26, 4).Style = "Normal"
$ws.Cells.Item(26, 5).Value = "  -1.90%  "
$ws.Cells.Item(27, 5).Value = "  -6.07%  "
$ws.Cells.Item(28, 5).Value = "  +0.04%  "
$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = "170.32"
$ws.Cells.Item(29, 4).Style = "Normal"
$ws.Cells.Item(29, 5).Value = "  -0.28%  "
$ws.Cells.Item(30, 4).Value = "0.0₃0734"
$ws.Cells.Item(30, 5).Value = "  -0.42%  "
$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = "6.13"
$ws.Cells.Item(31, 4).Style = "Normal"
$ws.Cells.Item(31, 5).Value = "  -0.81%  "
$ws.Cells.Item(32, 5).Value = "  -1.46%  "
$ws.Cells.Item(33, 5).Value = "  -1.20%  "
$ws.Cells.Item(34, 5).Value = "  -0.05%  "
$ws.Cells.Item(35, 5).Value = "  +0.47%  "
$ws.Cells.Item(36, 5).Value = "  +0.54%  "
$ws.Cells.Item(37, 5).Value = "  -2.42%  "
$ws.Cells.Item(38, 5).Value = "  -0.73%  "
$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = "39.07"
$ws.Cells.Item(39, 4).Style = "Normal"
$ws.Cells.Item(39, 5).Value = "  -0.89%  "
$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = "141.26"
$ws.Cells.Item(40, 4).Style = "Normal"
$ws.Cells.Item(40, 5).Value = "  -5.92%  "
$ws.Cells.Item(41, 5).Value = "  -0.63%  "
$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = "288.67"
$ws.Cells.Item(42, 4).Style = "Normal"
$ws.Cells.Item(42, 5).Value = "  +1.55%  "
$ws.Cells.Item(43, 5).Value = "  +0.19%  "
$ws.Cells.Item(44, 5).Value = "  +1.20%  "
$ws.Cells.Item(45, 5).Value = "  +0.64%  "
$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = "18.96"
$ws.Cells.Item(46, 4).Style = "Normal"
$ws.Cells.Item(46, 5).Value = "  -1.96%  "
$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = "0.567"
$ws.Cells.Item(47, 4).Style = "Normal"
$ws.Cells.Item(47, 5).Value = "  +0.62%  "
$ws.Cells.Item(48, 5).Value = "  +1.59%  "
$ws.Cells.Item(49, 5).Value = "  +0.11%  "
$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = "11.07"
$ws.Cells.Item(50, 4).Style = "Normal"
$ws.Cells.Item(50, 5).Value = "  +0.18%  "
$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = "0.954"
$ws.Cells.Item(51, 4).Style = "Normal"
$ws.Cells.Item(51, 5).Value = "  +0.84%  "
